$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell J1 ("Q8") - copy formatting from I1 (bold/border/center) then set value
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("J1").Value = "Q8"

# Updated / new numeric cell values (simulated rt_data bugfix)
$ws.Range("B2").Value = -0.8350711388588363
$ws.Range("C2").Value = 0.2201878744147194
$ws.Range("D2").Value = 0.4384779472729292
$ws.Range("E2").Value = 1.529210744895579
$ws.Range("F2").Value = 1.94593406307078
$ws.Range("G2").Value = 0.3611512668594514
$ws.Range("H2").Value = 1.308369230529848
$ws.Range("B3").Value = -0.0806700814468968
$ws.Range("C3").Value = 0.137619991411313
$ws.Range("D3").Value = 1.228352789033963
$ws.Range("E3").Value = 1.645076107209164
$ws.Range("F3").Value = 0.06029331099783519
$ws.Range("G3").Value = 1.007511274668232
$ws.Range("B4").Value = -0.43807954759518
$ws.Range("C4").Value = 0.6526532500274698
$ws.Range("D4").Value = 1.069376568202671
$ws.Range("E4").Value = -0.5154062280086578
$ws.Range("F4").Value = 0.4318117356617392
$ws.Range("G4").Value = 0.3521176394643107
$ws.Range("H4").Value = 0.06357296580725347
$ws.Range("I4").Value = 0.05403164090613188
$ws.Range("J4").Value = -0.9652956092697305
$ws.Range("B5").Value = -0.4230865868247518
$ws.Range("C5").Value = -0.006363268649550946
$ws.Range("D5").Value = -1.591146064860879
$ws.Range("E5").Value = -0.6439281011904823
$ws.Range("F5").Value = -0.7236221973879109
$ws.Range("G5").Value = -1.012166871044968
$ws.Range("H5").Value = -1.02170819594609
$ws.Range("I5").Value = -2.041035446121952
$ws.Range("B6").Value = 0.4354409026540649
$ws.Range("C6").Value = -1.149341893557263
$ws.Range("D6").Value = -0.2021239298868664
$ws.Range("E6").Value = -0.281818026084295
$ws.Range("F6").Value = -0.5703626997413522
$ws.Range("G6").Value = -0.5799040246424738
$ws.Range("H6").Value = -1.599231274818336
$ws.Range("B7").Value = -0.3395690612336324
$ws.Range("C7").Value = 0.6076489024367646
$ws.Range("D7").Value = 0.527954806239336
$ws.Range("E7").Value = 0.2394101325822788
$ws.Range("F7").Value = 0.2298688076811572
$ws.Range("G7").Value = -0.7894584424947051
$ws.Range("B8").Value = 0.3275174550509519
$ws.Range("C8").Value = 0.2478233588535233
$ws.Range("D8").Value = -0.04072131480353391
$ws.Range("E8").Value = -0.05026263970465549
$ws.Range("F8").Value = -1.069589889880518
$ws.Range("G8").Value = -1.047781072492029
$ws.Range("H8").Value = 0.4911891579199903
$ws.Range("I8").Value = -0.4608111822815033
$ws.Range("B9").Value = 0.2089562936243113
$ws.Range("C9").Value = -0.07958838003274593
$ws.Range("D9").Value = -0.08912970493386752
$ws.Range("E9").Value = -1.10845695510973
$ws.Range("F9").Value = -1.086648137721241
$ws.Range("G9").Value = 0.4523220926907783
$ws.Range("H9").Value = -0.4996782475107153
$ws.Range("B10").Value = 0.02912383308249389
$ws.Range("C10").Value = 0.0195825081813723
$ws.Range("D10").Value = -0.99974474199449
$ws.Range("E10").Value = -0.9779359246060011
$ws.Range("F10").Value = 0.5610343058060181
$ws.Range("G10").Value = -0.3909660343954755
$ws.Range("B11").Value = -0.1406678742931149
$ws.Range("C11").Value = -1.159995124468977
$ws.Range("D11").Value = -1.138186307080488
$ws.Range("E11").Value = 0.4007839233315309
$ws.Range("F11").Value = -0.5512164168699627
$ws.Range("B12").Value = -1.067943258854512
$ws.Range("C12").Value = -1.046134441466023
$ws.Range("D12").Value = 0.492835788945996
$ws.Range("E12").Value = -0.4591645512554976
$ws.Range("B13").Value = -0.8299273031874748
$ws.Range("C13").Value = 0.7090429272245444
$ws.Range("D13").Value = -0.2429574129769492
$ws.Range("B14").Value = 1.026563613664763
$ws.Range("C14").Value = 0.0745632734632693
$ws.Range("B15").Value = -0.2176617297482864
